# "Begin separating Control Layer" - minor text cleanups and date bump.
$d = $word.ActiveDocument

# 1) Merge the split runs in the "Target Audience" line (no textual change,
#    just forces Word to recombine adjacent identically-formatted runs).
$d.Content.Find.Execute(
    " 13 yrs.+ fans of old school RPGs like Diablo. People who like comedy. Rated T for Teens.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " 13 yrs.+ fans of old school RPGs like Diablo. People who like comedy. Rated T for Teens.",
    2) | Out-Null

# 2) Merge the split runs making up the first part of the Game Summary
#    paragraph (up to the italic "Wizard magic" phrase), and the italic
#    "Wizard magic" phrase itself (same paragraph, one pass covers both).
$d.Content.Find.Execute(
    " Project “Dragon” is a Diablo style RPG with the absurd humor of Monty Python and the Holy Grail. The game is set during the grubby phase of the medieval era when brave Sir Knights were bold and foolish, and the unbathed peasants reeked of cabbage. As our unlikely, elderly hero, you must battle the Queen’s soldiers and grotesque monsters by finding ancient weapons, leveling up your character, venturing to far off lands (using special “",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Project “Dragon” is a Diablo style RPG with the absurd humor of Monty Python and the Holy Grail. The game is set during the grubby phase of the medieval era when brave Sir Knights were bold and foolish, and the unbathed peasants reeked of cabbage. As our unlikely, elderly hero, you must battle the Queen’s soldiers and grotesque monsters by finding ancient weapons, leveling up your character, venturing to far off lands (using special “",
    2) | Out-Null

# 3) Merge the split runs making up the italic "Wizard magic" phrase.
$d.Content.Find.Execute(
    "Wizard magic” ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Wizard magic” ",
    2) | Out-Null

# 4) Merge the split runs in the "Anticipated Remarkability" line.
$d.Content.Find.Execute(
    " Absurd humor.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Absurd humor.",
    2) | Out-Null

# 5) Bump the anticipated launch year -- but keep it as its own run,
#    separate from the preceding "End of " run, by nudging a character
#    format after the text swap so the editor doesn't coalesce the two
#    adjacent runs back together.
$yearRng = $d.Content
$yearRng.Find.Execute("2019") | Out-Null
$yearRng.Text = "2020"
$origColor = $yearRng.Font.Color
$yearRng.Font.Color = 255
$yearRng.Font.Color = $origColor

# 6) Drop the stray _GoBack bookmark left over from the last edit session.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
